$d = $word.ActiveDocument

# 1. Title change
$d.Content.Find.Execute(
    "Too Many Threes? A Data-Driven NBA Debate", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The Three-Point Revolution: Is It Changing Basketball for Better or Worse?", 2)

# 2. NBA Statistics bullet - collapse the split runs into one run (text unchanged)
$d.Content.Find.Execute(
    "NBA Statistics: Three-point attempts per game, mid-range shot frequency, shooting percentages, game outcomes (score differentials, blowouts, overtime games, competitiveness), and player shot selection trends; sources: Basketball Reference, NBA statistics",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "NBA Statistics: Three-point attempts per game, mid-range shot frequency, shooting percentages, game outcomes (score differentials, blowouts, overtime games, competitiveness), and player shot selection trends; sources: Basketball Reference, NBA statistics",
    2)

# 3. TV Ratings & Fan Engagement bullet - collapse the split runs after the heading
#    (keep the lastRenderedPageBreak-holding run separate so it is untouched)
$d.Content.Find.Execute(
    " Yearly NBA Finals, regular-season, and playoff ratings, attendance data and ticket sales over time; sources: Sports Business Journal, ESPN, NBA revenue reports",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Yearly NBA Finals, regular-season, and playoff ratings, attendance data and ticket sales over time; sources: Sports Business Journal, ESPN, NBA revenue reports",
    2)

# 4. Social Media Sentiment Data bullet - collapse runs on either side of the
#    proofErr markers (the proofErr-wrapped "Social media" run is untouched)
$d.Content.Find.Execute(
    "Social Media Sentiment Data: ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Social Media Sentiment Data: ", 2)

$d.Content.Find.Execute(
    " and Reddit discussions on three-point shooting, game excitement, and competitiveness, text analysis of fan opinions on whether the game has become less entertaining; sources: Reddit, YouTube comments on NBA debates, other social media platforms.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " and Reddit discussions on three-point shooting, game excitement, and competitiveness, text analysis of fan opinions on whether the game has become less entertaining; sources: Reddit, YouTube comments on NBA debates, other social media platforms.",
    2)

# 5. "I expect the structured data..." paragraph - collapse the split runs
$d.Content.Find.Execute(
    "I expect the structured data (NBA stats, TV ratings) to be a smaller dataset, while the social media data could be 500MB+ depending on the volume of scraped text. Data collection and preprocessing should take 1-2 weeks.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I expect the structured data (NBA stats, TV ratings) to be a smaller dataset, while the social media data could be 500MB+ depending on the volume of scraped text. Data collection and preprocessing should take 1-2 weeks.",
    2)

# 6. Timeline table - "Week 1" deliverable cell
$d.Content.Find.Execute(
    "Project proposal and research plan.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Project proposal and research plan.", 2)

# 7. Timeline table - "Week 2" deliverable cell
$d.Content.Find.Execute(
    ". Clean and preprocess data (handle missing values, structure datasets). Start exploratory data analysis (EDA).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ". Clean and preprocess data (handle missing values, structure datasets). Start exploratory data analysis (EDA).",
    2)

# 8. Timeline table - "Week 3" deliverable cell
$d.Content.Find.Execute(
    "Continue EDA, create initial visualizations of three-point trends.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Continue EDA, create initial visualizations of three-point trends.", 2)
